# Generate Report for Handback
# Updates timestamps/status recorded in the handback status workbook.
# Note: several cells across sheets shared the same underlying cached
# string value in the source workbook, so a single logical edit shows
# up in more than one place.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$overview = $wb.Sheets.Item("Overview")
# "Latest HO Xliff Generate Date" for the 1add3613... and c42c1da2... rows
$overview.Range("G2").Value = "2016-08-18 16:13:57"
$overview.Range("G4").Value = "2016-08-18 16:13:57"

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Sheets.Item("zh-cn")
# Priority changed from "ht" to "mt"
$zhcn.Range("E2").Value = "mt"
$zhcn.Range("E4").Value = "mt"
# Correspond Handoff Datetime
$zhcn.Range("H2").Value = "2016-08-18 16:13:52"
$zhcn.Range("H4").Value = "2016-08-18 16:13:52"
# Correspond Handback DateTime
$zhcn.Range("K2").Value = "2016-08-18 16:14:17"
$zhcn.Range("K4").Value = "2016-08-18 16:14:17"

# --- de-de sheet ------------------------------------------------------
$dede = $wb.Sheets.Item("de-de")
# Priority changed from "ht" to "mt" (same cached value as zh-cn)
$dede.Range("E2").Value = "mt"
$dede.Range("E4").Value = "mt"
# Correspond Handoff Datetime (same cached value as Overview's
# "Latest HO Xliff Generate Date")
$dede.Range("H2").Value = "2016-08-18 16:13:57"
$dede.Range("H4").Value = "2016-08-18 16:13:57"
# Correspond Handback DateTime
$dede.Range("K2").Value = "2016-08-18 16:14:24"
$dede.Range("K4").Value = "2016-08-18 16:14:24"
